$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_MB")
$ws.Name = "CRF_MB"
